# Adding problem 2 part 3 (draft of visualization)
#
# - Rename "Sheet1" to "personal_data_sheet"
# - Minimize the workbook window
# - Move the active selection to A2 (away from the previous I39)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet
$ws.Name = "personal_data_sheet"

# Minimize the workbook window
$excel.WindowState = -4140  # xlMinimized

# Make sure the sheet is active, then move the selection to A2
$ws.Activate()
$ws.Range("A2").Select()
